# Inital draft for part 1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update SMS_received (column M) values for rows 2-4 from 0 to 1
$ws.Range("M2").Value = 1
$ws.Range("M3").Value = 1
$ws.Range("M4").Value = 1

# Update No-show (column N) values for rows 5-7 from "No" to "Yes"
$ws.Range("N5").Value = "Yes"
$ws.Range("N6").Value = "Yes"
$ws.Range("N7").Value = "Yes"

# Move the active cell selection to M9
$ws.Range("M9").Select()
